{"js": "// Find the final paragraph in the document body (\"Design decisions are\n// variable, ...\"), drop its trailing full stop, then insert a large block\n// of additional meeting notes (new paragraphs) immediately after it, still\n// before the end of the document / section break.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items[paragraphs.items.length - 1];\nconst currentText = target.text.replace(/[\\r\\u0007]+$/, \"\");\n\nif (!currentText.endsWith(\n  \"Design decisions are variable, but how you justify it or explain it is where the marks come from.\"\n)) {\n  throw new Error(\"Unexpected final paragraph; aborting to avoid editing the wrong content.\");\n}\n\n// Strip the trailing \".\" while leaving the rest of the run/paragraph as-is.\nconst trimmedText = currentText.slice(0, -1);\ntarget.getRange().insertText(trimmedText, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-fetch the (now trimmed) paragraph's end range and append the new\n// OOXML content \u2014 a flat-OPC WordprocessingML package \u2014 right after it.\nconst endRange = target.getRange(\"End\");\nconst xml = \"<?xml version=\\\"1.0\\\" encoding=\\\"UTF-8\\\" standalone=\\\"yes\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p/><w:p><w:pPr><w:jc w:val=\\\"center\\\"/><w:rPr><w:b/><w:bCs/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\\\"28\\\"/><w:szCs w:val=\\\"28\\\"/></w:rPr><w:t>27/09/2022</w:t></w:r></w:p><w:p><w:r><w:t>Focus is on package holdiers \\u201cglorified hotel database\\u201d</w:t></w:r></w:p><w:p><w:r><w:t>If a room is available assume corresponding flight</w:t></w:r></w:p><w:p><w:r><w:t>Things like guestTypes tables with only a few rows is fine</w:t></w:r></w:p><w:p><w:r><w:t>Play with stuff, eg prices, check with different numbers of adults and kids etc, see how it affects the price.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:t>would you take into account different prices from different airports? Or treat all airports the same?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t>Is it the date that effects it, is it the distance, is this stored somewhere as prices for each permutation, or is it calculated?</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:t>should we be factoring in different prices at different times of the year? and include flights with two legs or just assume direct flights?</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> Reasonable assumptions about to leave out. Don\\u2019t touch multi leg flights</w:t></w:r><w:r><w:t>, save your sanity.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> Different time of year is worth looking into.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:lastRenderedPageBreak/><w:t>so in a pricing table should you be doing like a base price for the package holiday then price adjustments calculated based on % at different ties of years and different airports?</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> Storing multipliers and stuff?</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> Acceptable. Demonstrating dynamic pricing based on dates.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Table numbers have often ended up in the range of 20-30. This is not a hard rule, and some of these tables are small, some are large etc</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:t>2 questions - from memory on Jet2 you are able to save multiple holidays for later without buying - should we facilitate this?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t>Yes and no, its viable information, but do you have time to speak about it in 10 minute video?</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:t>Also you can pay a deposit &amp; settle balance later - should we allow for this?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t>Sure, it may not even be that difficult to implement Booksing contain cost, and payments are made in payments table.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> Monthly payment algorithms craic not in scope.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:t>Do we need to worry about concurrent bookings with limited availability (like ticketmaster reserving a ticket and giving you 15 minutes to check out)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Segoe UI\\\" w:hAnsi=\\\"Segoe UI\\\" w:cs=\\\"Segoe UI\\\"/><w:color w:val=\\\"FFFFFF\\\"/><w:sz w:val=\\\"21\\\"/><w:szCs w:val=\\\"21\\\"/><w:shd w:val=\\\"clear\\\" w:color=\\\"auto\\\" w:fill=\\\"292929\\\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> Ticketmaster has presumptive booking.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> A selected holiday on website may store a in progress booking with a status saying like \\u201cunconfirmed\\u201d. Payment just switches status from unconfirmed to confirmed</w:t></w:r><w:r><w:t>. There is a concurrency management system in place, but the limitations of it are largely unknown from the outside</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Room Types: are the particular to hotels, or are they hard categories across all hotels?</w:t></w:r></w:p><w:p><w:r><w:t>Forget about luggage and sports equipment.</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\nendRange.insertOoxml(xml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Locate the final paragraph (\"Design decisions are variable, ...\") and\n# strip its trailing full stop, then append a large block of new meeting\n# notes (as a fresh OOXML fragment) immediately after it, before the\n# section break.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Design decisions are variable*marks come from*\") {\n        $target = $p\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Design decisions...' paragraph\"\n}\n\n$r = $target.Range\n# Trim the trailing paragraph mark off the range so Find/char ops only touch\n# the visible text, then drop the final \".\" character.\n$textRange = $r.Duplicate\n$textRange.MoveEnd(1, -1) | Out-Null\nif ($textRange.Text.EndsWith(\".\")) {\n    $textRange.MoveEnd(1, -1) | Out-Null\n    $charRange = $d.Range($textRange.End, $textRange.End + 1)\n    $charRange.Delete() | Out-Null\n}\n\n$insertPos = $d.Content.End\n$insertRange = $d.Range($insertPos, $insertPos)\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p/><w:p><w:pPr><w:jc w:val=\"center\"/><w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>27/09/2022</w:t></w:r></w:p><w:p><w:r><w:t>Focus is on package holdiers \u201cglorified hotel database\u201d</w:t></w:r></w:p><w:p><w:r><w:t>If a room is available assume corresponding flight</w:t></w:r></w:p><w:p><w:r><w:t>Things like guestTypes tables with only a few rows is fine</w:t></w:r></w:p><w:p><w:r><w:t>Play with stuff, eg prices, check with different numbers of adults and kids etc, see how it affects the price.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:t>would you take into account different prices from different airports? Or treat all airports the same?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>Is it the date that effects it, is it the distance, is this stored somewhere as prices for each permutation, or is it calculated?</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:t>should we be factoring in different prices at different times of the year? and include flights with two legs or just assume direct flights?</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Reasonable assumptions about to leave out. Don\u2019t touch multi leg flights</w:t></w:r><w:r><w:t>, save your sanity.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Different time of year is worth looking into.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:lastRenderedPageBreak/><w:t>so in a pricing table should you be doing like a base price for the package holiday then price adjustments calculated based on % at different ties of years and different airports?</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Storing multipliers and stuff?</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Acceptable. Demonstrating dynamic pricing based on dates.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Table numbers have often ended up in the range of 20-30. This is not a hard rule, and some of these tables are small, some are large etc</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:t>2 questions - from memory on Jet2 you are able to save multiple holidays for later without buying - should we facilitate this?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>Yes and no, its viable information, but do you have time to speak about it in 10 minute video?</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:t>Also you can pay a deposit &amp; settle balance later - should we allow for this?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>Sure, it may not even be that difficult to implement Booksing contain cost, and payments are made in payments table.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Monthly payment algorithms craic not in scope.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:t>Do we need to worry about concurrent bookings with limited availability (like ticketmaster reserving a ticket and giving you 15 minutes to check out)</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Segoe UI\" w:hAnsi=\"Segoe UI\" w:cs=\"Segoe UI\"/><w:color w:val=\"FFFFFF\"/><w:sz w:val=\"21\"/><w:szCs w:val=\"21\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"292929\"/></w:rPr><w:t>.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Ticketmaster has presumptive booking.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> A selected holiday on website may store a in progress booking with a status saying like \u201cunconfirmed\u201d. Payment just switches status from unconfirmed to confirmed</w:t></w:r><w:r><w:t>. There is a concurrency management system in place, but the limitations of it are largely unknown from the outside</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Room Types: are the particular to hotels, or are they hard categories across all hotels?</w:t></w:r></w:p><w:p><w:r><w:t>Forget about luggage and sports equipment.</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$insertRange.InsertXML($xml) | Out-Null\n"}
